# Updates after conversation with Keith
#
# Inserts a new leading "Element ID" column into the domain model table,
# shifting the existing Name/Owner/Base Classifier/Type/Documentation
# columns one position to the right, then populates the new column with
# per-row element identifiers. Also refreshes a couple of relationship
# rows whose Owner/Base Classifier values changed alongside the shift.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new column at A, pushing the old A:E data to B:F.
$ws.Columns.Item(1).Insert()

# 2. Give the new column the same formatting (wrap text, left/top align)
#    as the data that just got shifted into column B.
$ws.Range("B1:B44").Copy()
$ws.Range("A1:A44").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# 3. Populate the new "Element ID" column.
$ws.Range("A1").Value = "Element ID"
$ws.Range("A2").Value = "_19_0_3_62501eb_1575932667283_578110_6044"
$ws.Range("A3").Value = "_19_0_3_62501eb_1575932667284_121545_6045"
$ws.Range("A4").Value = "_19_0_3_62501eb_1578694996433_774413_6726"
$ws.Range("A5").Value = "_19_0_3_62501eb_1578694996434_437053_6727"
$ws.Range("A6").Value = "_19_0_3_62501eb_1581105741408_455163_7331"
$ws.Range("A7").Value = "_19_0_3_62501eb_1576171289956_933528_12331"
$ws.Range("A8").Value = "_19_0_3_62501eb_1575932667401_681378_6304"
$ws.Range("A9").Value = "_19_0_3_62501eb_1581439410653_635908_6762"
$ws.Range("A10").Value = "_19_0_3_62501eb_1581439888040_737322_6776"
$ws.Range("A11").Value = "_19_0_3_62501eb_1575938781064_682315_11139"
$ws.Range("A12").Value = "_19_0_3_62501eb_1581442515746_351649_6834"
$ws.Range("A13").Value = "_19_0_3_62501eb_1575932667280_478117_6038"
$ws.Range("A14").Value = "_19_0_3_62501eb_1575932667280_396639_6039"
$ws.Range("A15").Value = "_19_0_3_62501eb_1581439410653_264931_6761"
$ws.Range("A16").Value = "_19_0_3_62501eb_1581442515746_256830_6833"
$ws.Range("A17").Value = "_19_0_3_62501eb_1575938781063_864467_11138"
$ws.Range("A18").Value = "_19_0_3_62501eb_1581115363833_428095_6704"
$ws.Range("A19").Value = "_19_0_3_62501eb_1581105619554_611153_7311"
$ws.Range("A20").Value = "_19_0_3_62501eb_1581105562224_754275_7303"
$ws.Range("A21").Value = "_19_0_3_62501eb_1581105916125_511300_7348"
$ws.Range("A22").Value = "_19_0_3_62501eb_1575932667400_432786_6303"
$ws.Range("A23").Value = "_19_0_3_62501eb_1581105641554_143458_7315"
$ws.Range("A24").Value = "_19_0_3_62501eb_1576171289955_293204_12330"
$ws.Range("A25").Value = "_19_0_3_62501eb_1581105588214_737974_7307"
$ws.Range("A26").Value = "_19_0_3_62501eb_1576170905738_395553_12292"
$ws.Range("A27").Value = "_19_0_3_62501eb_1581116150135_739112_6711"
$ws.Range("A28").Value = "_19_0_3_62501eb_1581105741408_741963_7330"
$ws.Range("A29").Value = "_19_0_3_62501eb_1575932666937_126818_5370"
$ws.Range("A30").Value = "_19_0_3_62501eb_1575932666914_330525_5327"
$ws.Range("A31").Value = "_19_0_3_62501eb_1578694912303_354692_6680"
$ws.Range("A32").Value = "_19_0_3_62501eb_1578694889873_510687_6655"
$ws.Range("A33").Value = "_19_0_3_62501eb_1575932666853_60531_5213"
$ws.Range("A34").Value = "_19_0_3_62501eb_1576709459751_212403_6610"
$ws.Range("A35").Value = "_19_0_3_62501eb_1575932666803_646860_5127"
$ws.Range("A36").Value = "_19_0_3_62501eb_1576848426080_581998_6606"
$ws.Range("A37").Value = "_19_0_3_62501eb_1575932666920_674563_5339"
$ws.Range("A38").Value = "_19_0_3_62501eb_1576101195264_864774_11895"
$ws.Range("A39").Value = "_19_0_3_62501eb_1575932666733_316702_5026"
$ws.Range("A40").Value = "_19_0_3_62501eb_1575938749735_717594_11092"
$ws.Range("A41").Value = "_19_0_3_62501eb_1575932666846_481711_5201"
$ws.Range("A42").Value = "_19_0_3_62501eb_1576709434181_966028_6577"
$ws.Range("A43").Value = "_19_0_3_62501eb_1575932666765_847074_5072"
$ws.Range("A44").Value = "_19_0_3_62501eb_1575932666745_678311_5042"

# 4. A couple of rows also had their Owner / Base Classifier content
#    revised (not just shifted) during the "conversation with Keith"
#    clean-up pass.

# Row 15: "vital sign" used to be subject's owner here; it is now owned
# by "health related condition" instead, and the old "condition of"
# relationship name (row 16) is dropped.
$ws.Range("C15").Value = "health related condition"
$ws.Range("C16").Value = "vital sign"
$ws.Range("B16").Value = ""

# Row 43: "vital sign" now additionally classifies as
# "health related condition" alongside "quality".
$ws.Range("D43").Value = "quality`nhealth related condition"
